$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.724.86'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '3.446.73'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.21'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.97'
$ws.Range("E6").Value = '  +1.57%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.487'
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.03'
$ws.Range("E9").Value = '  +5.55%  '
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("E11").Value = '  +3.73%  '
$ws.Range("D12").Value = '4.038.53'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.18'
$ws.Range("E14").Value = '  -5.17%  '
$ws.Range("D15").Value = '3.444.84'
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '62.746.73'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.38'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.62'
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.97'
$ws.Range("E20").Value = '  -2.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '387.39'
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.567'
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = '3.586.87'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.182'
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  -4.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.24'
$ws.Range("E34").Value = '  -2.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.42'
$ws.Range("E35").Value = '  +2.01%  '
$ws.Range("E36").Value = '  +3.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '31.51'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.93'
$ws.Range("E38").Value = '  -2.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '168.95'
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("D40").Value = '3.482.41'
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0777'
$ws.Range("E41").Value = '  +1.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.85'
$ws.Range("E42").Value = '  +1.62%  '
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("E45").Value = '  -1.55%  '
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("D47").Value = '2.564.08'
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("E48").Value = '  +2.22%  '
$ws.Range("E49").Value = '  +1.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.69'
$ws.Range("E50").Value = '  -3.07%  '
$ws.Range("E51").Value = '  +0.06%  '
